# Update crypto symbol list (scraped values refreshed by the GitHub Actions
# job on Fri Jan 6 20:47:20 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume(1h) columns hold numeric- and percent-looking text that
# Excel would otherwise auto-convert into real numbers on assignment. Mark
# those specific cells as Text first so the literal strings from the source
# scrape survive, exactly as they were originally stored.
$numericLookingCells = @("D2","E2","E3","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","E26","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50")
foreach ($nc in $numericLookingCells) {
    $ws.Range($nc).NumberFormat = "@"
}

$updates = @(
    @{ Cell = "D2";  Value = "259.22" },
    @{ Cell = "E2";  Value = "0.60%" },
    @{ Cell = "E3";  Value = "-0.89%" },
    @{ Cell = "E4";  Value = "0.09%" },
    @{ Cell = "D5";  Value = "0.06035" },
    @{ Cell = "E5";  Value = "2.42%" },
    @{ Cell = "D6";  Value = "6.681" },
    @{ Cell = "E6";  Value = "0.58%" },
    @{ Cell = "D7";  Value = "0.8583" },
    @{ Cell = "E7";  Value = "0.04%" },
    @{ Cell = "D8";  Value = "0.9305" },
    @{ Cell = "E8";  Value = "-1.43%" },
    @{ Cell = "D9";  Value = "0.1394" },
    @{ Cell = "E9";  Value = "-1.00%" },
    @{ Cell = "D10"; Value = "0.04762" },
    @{ Cell = "E10"; Value = "21.67%" },
    @{ Cell = "D11"; Value = "0.07082" },
    @{ Cell = "E11"; Value = "-0.03%" },
    @{ Cell = "D12"; Value = "0.03153" },
    @{ Cell = "E12"; Value = "-0.72%" },
    @{ Cell = "D13"; Value = "0.09131" },
    @{ Cell = "E13"; Value = "-0.34%" },
    @{ Cell = "D14"; Value = "0.001534" },
    @{ Cell = "E14"; Value = "-0.92%" },
    @{ Cell = "D15"; Value = "0.0006057" },
    @{ Cell = "E15"; Value = "0.49%" },
    @{ Cell = "D16"; Value = "0.006008" },
    @{ Cell = "E16"; Value = "-3.13%" },
    @{ Cell = "D17"; Value = "3.463" },
    @{ Cell = "E17"; Value = "-1.40%" },
    @{ Cell = "D18"; Value = "3.165" },
    @{ Cell = "E18"; Value = "-1.15%" },
    @{ Cell = "D19"; Value = "2.186" },
    @{ Cell = "E19"; Value = "-1.74%" },
    @{ Cell = "D21"; Value = "0.1298" },
    @{ Cell = "E21"; Value = "0.42%" },
    @{ Cell = "D22"; Value = "4.120" },
    @{ Cell = "E22"; Value = "6.22%" },
    @{ Cell = "D23"; Value = "0.04238" },
    @{ Cell = "E23"; Value = "0.30%" },
    @{ Cell = "D24"; Value = "0.001217" },
    @{ Cell = "E24"; Value = "-0.43%" },
    @{ Cell = "D25"; Value = "0.004040" },
    @{ Cell = "E25"; Value = "-5.99%" },
    @{ Cell = "E26"; Value = "-0.03%" },
    @{ Cell = "D40"; Value = "0.03843" },
    @{ Cell = "E40"; Value = "0.11%" },
    @{ Cell = "D41"; Value = "0.1117" },
    @{ Cell = "E41"; Value = "1.33%" },
    @{ Cell = "D42"; Value = "0.003939" },
    @{ Cell = "E42"; Value = "-36.86%" },
    @{ Cell = "B43"; Value = "LocalTraders" },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct" },
    @{ Cell = "D43"; Value = "0.01529" },
    @{ Cell = "E43"; Value = "33.37%" },
    @{ Cell = "B44"; Value = "CEJI" },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji" },
    @{ Cell = "D44"; Value = "0.002199" },
    @{ Cell = "E44"; Value = "-0.01%" },
    @{ Cell = "E45"; Value = "-6.51%" },
    @{ Cell = "D46"; Value = "0.00000000750" },
    @{ Cell = "E46"; Value = "-0.01%" },
    @{ Cell = "D47"; Value = "0.05454" },
    @{ Cell = "E47"; Value = "-9.06%" },
    @{ Cell = "D48"; Value = "0.1321" },
    @{ Cell = "E48"; Value = "2.66%" },
    @{ Cell = "D49"; Value = "0.00002099" },
    @{ Cell = "E49"; Value = "-0.01%" },
    @{ Cell = "D50"; Value = "0.0001999" },
    @{ Cell = "E50"; Value = "-0.01%" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
